$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry the per-observation data which gets rotated
# across rows 30-35. L and N are always-empty placeholder cells that never
# change, so they are intentionally left untouched.
$cols = "A","B","E","F","G","H","I","J","K","M","Q","R"

# Capture current values (and "does the cell have content" flags) for rows
# 30-35 across the moved columns before any writes happen.
$data = @{}
foreach ($r in 30..35) {
    $rowData = @{}
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $rowData[$col] = $cell.Value2
    }
    $data[$r] = $rowData
}

# mapping: new row -> old row (where the content should come from)
$mapping = @{30=33; 31=35; 32=34; 33=32; 34=30; 35=31}

foreach ($newRow in 30..35) {
    $oldRow = $mapping[$newRow]
    $src = $data[$oldRow]
    foreach ($col in $cols) {
        $val = $src[$col]
        if ($null -eq $val) {
            $ws.Range("$col$newRow").ClearContents()
        } else {
            $ws.Range("$col$newRow").Value2 = $val
        }
    }
}
